# 20T to 21 T update fiscal year.
# The "pe" (period) column (M) on the "19Tto20TMap" sheet listed the fiscal
# year-start month "2018Oct" for every mapped indicator row; bump it forward
# one fiscal year to "2019Oct".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("19Tto20TMap")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 77 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 13)
    if ($cell.Value2 -eq "2018Oct") {
        $cell.Value = "2019Oct"
    }
}

# Restore the view state recorded in the saved workbook (best effort —
# scroll position is cosmetic only).
$ws.Activate()
$ws.Range("M66").Select()
